$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.554.06'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '2.929.66'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '350.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.38'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.554'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.81%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.603'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.01%  '
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0849'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.00%  '
$ws.Range("D14").Value = '3.386.26'
$ws.Range("E14").Value = '  -0.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.81%  '
$ws.Range("D16").Value = '2.917.67'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("E17").Value = '  -1.84%  '
$ws.Range("D18").Value = '51.527.99'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.81%  '
$ws.Range("E22").Value = '  -1.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '260.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = '  -2.75%  '
$ws.Range("E26").Value = '  -4.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.43'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.07%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.105'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '10.21'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.61%  '
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '35.59'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.42'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.73%  '
$ws.Range("E36").Value = '  -5.03%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("E38").Value = '  -7.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.67'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.61%  '
$ws.Range("E42").Value = '  -1.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.23'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.75'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.18%  '
$ws.Range("E45").Value = '  -3.27%  '
$ws.Range("D46").Value = '2.094.43'
$ws.Range("E46").Value = '  -4.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.30'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.29'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.237'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0335'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.903'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.94%  '
